$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.00841100000001
$ws.Range("H2").Value = 75.02523300000001
$ws.Range("I2").Value = 0.4156829172908309
$ws.Range("J2").Value = 0.415682917290831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 25.01567066666667
$ws.Range("N2").Value = 75.047012
$ws.Range("O2").Value = 0.3530689998156723
$ws.Range("P2").Value = 0.3530689998156723
$ws.Range("Q2").Value = 625.6021734726442
$ws.Range("R2").Value = 5630.419561253797
$ws.Range("S2").Value = 0.1467647518483345
$ws.Range("T2").Value = 0.1467647518483345

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.00841100000001
$ws.Range("H3").Value = 75.02523300000001
$ws.Range("I3").Value = 0.4156829172908309
$ws.Range("J3").Value = 0.415682917290831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 25.99788533333333
$ws.Range("N3").Value = 77.993656
$ws.Range("O3").Value = 0.3669318921836303
$ws.Range("P3").Value = 0.3669318921836303
$ws.Range("Q3").Value = 650.1658015468721
$ws.Range("R3").Value = 5851.492213921849
$ws.Range("S3").Value = 0.1525273193899361
$ws.Range("T3").Value = 0.1525273193899361

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.00841100000001
$ws.Range("H4").Value = 75.02523300000001
$ws.Range("I4").Value = 0.4156829172908309
$ws.Range("J4").Value = 0.415682917290831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.83851733333333
$ws.Range("N4").Value = 59.515552
$ws.Range("O4").Value = 0.2799991080006974
$ws.Range("P4").Value = 0.2799991080006974
$ws.Range("Q4").Value = 496.1297951026241
$ws.Range("R4").Value = 4465.168155923617
$ws.Range("S4").Value = 0.1163908460525603
$ws.Range("T4").Value = 0.1163908460525604

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.496019
$ws.Range("H5").Value = 61.488057
$ws.Range("I5").Value = 0.340679180727168
$ws.Range("J5").Value = 0.3406791807271681
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.01567066666667
$ws.Range("N5").Value = 75.047012
$ws.Range("O5").Value = 0.3530689998156723
$ws.Range("P5").Value = 0.3530689998156723
$ws.Range("Q5").Value = 512.7216612817426
$ws.Range("R5").Value = 4614.494951535684
$ws.Range("S5").Value = 0.1202832575973639
$ws.Range("T5").Value = 0.1202832575973639

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.496019
$ws.Range("H6").Value = 61.488057
$ws.Range("I6").Value = 0.340679180727168
$ws.Range("J6").Value = 0.3406791807271681
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 25.99788533333333
$ws.Range("N6").Value = 77.993656
$ws.Range("O6").Value = 0.3669318921836303
$ws.Range("P6").Value = 0.3669318921836303
$ws.Range("Q6").Value = 532.8531517518213
$ws.Range("R6").Value = 4795.678365766392
$ws.Range("S6").Value = 0.1250060564117887
$ws.Range("T6").Value = 0.1250060564117887

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.496019
$ws.Range("H7").Value = 61.488057
$ws.Range("I7").Value = 0.340679180727168
$ws.Range("J7").Value = 0.3406791807271681
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 19.83851733333333
$ws.Range("N7").Value = 59.515552
$ws.Range("O7").Value = 0.2799991080006974
$ws.Range("P7").Value = 0.2799991080006974
$ws.Range("Q7").Value = 406.6106281958293
$ws.Range("R7").Value = 3659.495653762464
$ws.Range("S7").Value = 0.09538986671801544
$ws.Range("T7").Value = 0.09538986671801546

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.65779933333333
$ws.Range("H8").Value = 43.973398
$ws.Range("I8").Value = 0.243637901982001
$ws.Range("J8").Value = 0.243637901982001
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.01567066666667
$ws.Range("N8").Value = 75.047012
$ws.Range("O8").Value = 0.3530689998156723
$ws.Range("P8").Value = 0.3530689998156723
$ws.Range("Q8").Value = 366.6746808207529
$ws.Range("R8").Value = 3300.072127386776
$ws.Range("S8").Value = 0.08602099036997389
$ws.Range("T8").Value = 0.0860209903699739

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.65779933333333
$ws.Range("H9").Value = 43.973398
$ws.Range("I9").Value = 0.243637901982001
$ws.Range("J9").Value = 0.243637901982001
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.99788533333333
$ws.Range("N9").Value = 77.993656
$ws.Range("O9").Value = 0.3669318921836303
$ws.Range("P9").Value = 0.3669318921836303
$ws.Range("Q9").Value = 381.0717863070098
$ws.Range("R9").Value = 3429.646076763088
$ws.Range("S9").Value = 0.08939851638190548
$ws.Range("T9").Value = 0.08939851638190549

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.65779933333333
$ws.Range("H10").Value = 43.973398
$ws.Range("I10").Value = 0.243637901982001
$ws.Range("J10").Value = 0.243637901982001
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.83851733333333
$ws.Range("N10").Value = 59.515552
$ws.Range("O10").Value = 0.2799991080006974
$ws.Range("P10").Value = 0.2799991080006974
$ws.Range("Q10").Value = 290.7890061428551
$ws.Range("R10").Value = 2617.101055285696
$ws.Range("S10").Value = 0.06821839523012163
$ws.Range("T10").Value = 0.06821839523012165
